$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# --- Copy formatting down from the last existing data row (29) ---
$ws.Range("A29:E29").Copy()
$ws.Range("A30:E32").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A29:E29").Copy()
$ws.Range("A33:E35").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 30: TopicTypeaheadValidationTest ---
$ws.Cells.Item(30, 1).Value = "TopicTypeaheadValidationTest"
$ws.Cells.Item(30, 2).Value = "TBD"
$ws.Cells.Item(30, 3).Value = "Verify that topic typeahead options should display while enter min 2 characters"
$ws.Cells.Item(30, 4).Value = "Y"
$ws.Cells.Item(30, 5).Value = "SKIP"

# --- Row 31: ProfileTabFocusTest ---
$ws.Cells.Item(31, 1).Value = "ProfileTabFocusTest"
$ws.Cells.Item(31, 2).Value = "TBD"
$ws.Cells.Item(31, 3).Value = "Verify that By default profile tab foucs should be on POST tab"
$ws.Cells.Item(31, 4).Value = "Y"
$ws.Cells.Item(31, 5).Value = "SKIP"

# --- Row 32: PostTabInfiniteScrollValidationTest ---
$ws.Cells.Item(32, 1).Value = "PostTabInfiniteScrollValidationTest"
$ws.Cells.Item(32, 2).Value = "TBD"
$ws.Cells.Item(32, 3).Value = "Verify that Post tab infinite scroll displaying the more available records"
$ws.Cells.Item(32, 4).Value = "Y"
$ws.Cells.Item(32, 5).Value = "SKIP"

# --- Row 33: CommentsTabInfiniteScrollValidationTest ---
$ws.Cells.Item(33, 1).Value = "CommentsTabInfiniteScrollValidationTest"
$ws.Cells.Item(33, 2).Value = "TBD"
$ws.Cells.Item(33, 3).Value = "Verify that Comments tab infinite scroll displaying the more available records"
$ws.Cells.Item(33, 4).Value = "Y"
$ws.Cells.Item(33, 5).Value = "PASS"

# --- Row 34: FollowersTabInfiniteScrollValidationTest ---
$ws.Cells.Item(34, 1).Value = "FollowersTabInfiniteScrollValidationTest"
$ws.Cells.Item(34, 2).Value = "TBD"
$ws.Cells.Item(34, 3).Value = "Verify that Followers tab infinite scroll displaying the more available records"
$ws.Cells.Item(34, 4).Value = "Y"
$ws.Cells.Item(34, 5).Value = "PASS"

# --- Row 35: FollowingTabInfiniteScrollValidationTest ---
$ws.Cells.Item(35, 1).Value = "FollowingTabInfiniteScrollValidationTest"
$ws.Cells.Item(35, 2).Value = "TBD"
$ws.Cells.Item(35, 3).Value = "Verify that Following tab infinite scroll displaying the more available records"
$ws.Cells.Item(35, 4).Value = "Y"
$ws.Cells.Item(35, 5).Value = "PASS"

# --- Update the view state to match: scrolled so row 11 / col B is top-left, ---
# --- selection on D2:D35 with active cell D2 ---
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D2:D35").Select()
